$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ActiveBeneficiaryData-Mo")
$ws2 = $wb.Worksheets.Item("Sheet1")

# Column AH (Management Info Manager Name) gets filled for rows 2 through 182
# (row 183 is intentionally left blank in AH, matching the source edit).
$rngAH = $ws1.Range("AH2:AH182")
$rngAH.Value = "Management Info Manager Name"
$rngAH.Font.Bold = $true
$rngAH.HorizontalAlignment = -4130
$rngAH.VerticalAlignment = -4130
$rngAH.WrapText = $true

# Column AI (Management Info Business Partner Name) gets filled for rows 2 through 183.
$rngAI = $ws1.Range("AI2:AI183")
$rngAI.Value = "Management Info Business Partner Name"
$rngAI.Font.Bold = $true
$rngAI.HorizontalAlignment = -4130
$rngAI.VerticalAlignment = -4130
$rngAI.WrapText = $true

# Restore Sheet1's (second tab) selection back to B1 before leaving it.
$ws2.Activate()
$ws2.Range("B1").Select()

# Finish with the data sheet active and AH3 selected, matching the saved view state.
$ws1.Activate()
$ws1.Range("AH3").Select()
